$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newText = "Fix it in Baseboard_v0109, please verify"

# H4 already uses the red-font style (s="9") that we want to reuse for H7/H9
$ws.Range("H4").Copy()
$ws.Range("H7").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H7").Value = $newText

$ws.Range("H4").Copy()
$ws.Range("H9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H9").Value = $newText

$excel.CutCopyMode = 0

# Update the view state: scroll so row 9 is the top-left visible row,
# and select H9 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$ws.Range("H9").Select()
